$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

for ($r = 447; $r -le 526; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $cell.Value = "ok"
}
